# Auto-generated Excel COM-interop script applying scheduled-runner updates
# to the Ultima_Profits "Leve Profit" sheets (one per crafting class).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 2971
$ws.Range("I18").Value = 1949.25
$ws.Range("J18").Value = 4333.3335
$ws.Range("K18").Value = 1949.25
$ws.Range("L18").Value = 4333.3335
$ws.Range("M18").Value = -1665.25
$ws.Range("N18").Value = -4901.3335
$ws.Range("H19").Value = 1734.9166
$ws.Range("I19").Value = 1373.421
$ws.Range("J19").Value = 2138.9412
$ws.Range("K19").Value = 1373.421
$ws.Range("L19").Value = 2138.9412
$ws.Range("M19").Value = -1198.421
$ws.Range("N19").Value = -2488.9412
$ws.Range("H28").Value = 1319.3334
$ws.Range("I28").Value = 530.1818
$ws.Range("K28").Value = 530.1818
$ws.Range("M28").Value = -45.18179999999995
$ws.Range("H74").Value = 3612.3547
$ws.Range("I74").Value = 3326.9443
$ws.Range("J74").Value = 4007.5386
$ws.Range("K74").Value = 3326.9443
$ws.Range("L74").Value = 4007.5386
$ws.Range("M74").Value = -2390.9443
$ws.Range("N74").Value = -5879.5386
$ws.Range("H77").Value = 3612.3547
$ws.Range("I77").Value = 3326.9443
$ws.Range("J77").Value = 4007.5386
$ws.Range("K77").Value = 16634.7215
$ws.Range("L77").Value = 20037.693
$ws.Range("M77").Value = -11954.7215
$ws.Range("N77").Value = -29397.693
$ws.Range("H106").Value = 2698.5
$ws.Range("I106").Value = 2497.5
$ws.Range("J106").Value = 2899.5
$ws.Range("K106").Value = 2497.5
$ws.Range("L106").Value = 2899.5
$ws.Range("M106").Value = -1866.5
$ws.Range("N106").Value = -4161.5
$ws.Range("H112").Value = 1704.5143
$ws.Range("I112").Value = 675
$ws.Range("J112").Value = 1837.3549
$ws.Range("K112").Value = 2025
$ws.Range("L112").Value = 5512.0647
$ws.Range("M112").Value = -917
$ws.Range("N112").Value = -7728.0647
$ws.Range("H113").Value = 2227.7778
$ws.Range("I113").Value = 1827.2727
$ws.Range("J113").Value = 2857.1428
$ws.Range("K113").Value = 1827.2727
$ws.Range("L113").Value = 2857.1428
$ws.Range("M113").Value = 1426.7273
$ws.Range("N113").Value = -9365.1428

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1981
$ws.Range("I45").Value = 1962.6666
$ws.Range("J45").Value = 2011.5555
$ws.Range("K45").Value = 1962.6666
$ws.Range("L45").Value = 2011.5555
$ws.Range("M45").Value = -1585.6666
$ws.Range("N45").Value = -2765.5555
$ws.Range("H61").Value = 55559596
$ws.Range("I61").Value = 55559596
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 55559596
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -55559384
$ws.Range("N61").ClearContents()
$ws.Range("H122").Value = 7968.4443
$ws.Range("I122").Value = 9246.714
$ws.Range("J122").Value = 3494.5
$ws.Range("K122").Value = 27740.142
$ws.Range("L122").Value = 10483.5
$ws.Range("M122").Value = -25290.142
$ws.Range("N122").Value = -15383.5
$ws.Range("H132").Value = 6099146
$ws.Range("I132").Value = 7577132.5
$ws.Range("J132").Value = 2452
$ws.Range("K132").Value = 22731397.5
$ws.Range("L132").Value = 7356
$ws.Range("M132").Value = -22728867.5
$ws.Range("N132").Value = -12416
$ws.Range("H136").Value = 55559596
$ws.Range("I136").Value = 55559596
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 166678788
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -166676238
$ws.Range("N136").ClearContents()
$ws.Range("H140").Value = 36680
$ws.Range("J140").Value = 36680
$ws.Range("L140").Value = 36680
$ws.Range("N140").Value = -47040

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 35717212
$ws.Range("I86").Value = 2776
$ws.Range("J86").Value = 100003200
$ws.Range("K86").Value = 2776
$ws.Range("L86").Value = 100003200
$ws.Range("M86").Value = -1653
$ws.Range("N86").Value = -100005446
$ws.Range("H89").Value = 35717212
$ws.Range("I89").Value = 2776
$ws.Range("J89").Value = 100003200
$ws.Range("K89").Value = 13880
$ws.Range("L89").Value = 500016000
$ws.Range("M89").Value = -8264
$ws.Range("N89").Value = -500027232
$ws.Range("H107").Value = 1642.2
$ws.Range("I107").Value = 1437
$ws.Range("K107").Value = 1437
$ws.Range("M107").Value = 483

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 1402
$ws.Range("I99").Value = 1194.7142
$ws.Range("J99").Value = 1764.75
$ws.Range("K99").Value = 1194.7142
$ws.Range("L99").Value = 1764.75
$ws.Range("M99").Value = 303.2858000000001
$ws.Range("N99").Value = -4760.75
$ws.Range("H112").Value = 0
$ws.Range("J112").Value = 0
$ws.Range("L112").Value = 0
$ws.Range("N112").ClearContents()
$ws.Range("H126").Value = 1402
$ws.Range("I126").Value = 1194.7142
$ws.Range("J126").Value = 1764.75
$ws.Range("K126").Value = 3584.1426
$ws.Range("L126").Value = 5294.25
$ws.Range("M126").Value = -1114.1426
$ws.Range("N126").Value = -10234.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 505.15585
$ws.Range("I113").Value = 464.85455
$ws.Range("J113").Value = 605.9091
$ws.Range("K113").Value = 1394.56365
$ws.Range("L113").Value = 1817.7273
$ws.Range("M113").Value = 775.4363499999999
$ws.Range("N113").Value = -6157.7273
$ws.Range("H122").Value = 1005.8333
$ws.Range("I122").Value = 1020.4
$ws.Range("J122").Value = 933
$ws.Range("K122").Value = 9183.6
$ws.Range("L122").Value = 8397
$ws.Range("M122").Value = -6733.6
$ws.Range("N122").Value = -13297
$ws.Range("H129").Value = 2758.9412
$ws.Range("I129").Value = 1183
$ws.Range("J129").Value = 3167.5186
$ws.Range("K129").Value = 3549
$ws.Range("L129").Value = 9502.5558
$ws.Range("M129").Value = 1451
$ws.Range("N129").Value = -19502.5558
$ws.Range("H131").Value = 1239.6364
$ws.Range("I131").Value = 503.33334
$ws.Range("J131").Value = 1445.1163
$ws.Range("K131").Value = 1510.00002
$ws.Range("L131").Value = 4335.3489
$ws.Range("M131").Value = 3529.99998
$ws.Range("N131").Value = -14415.3489

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H19").Value = 8500
$ws.Range("J19").Value = 8500
$ws.Range("L19").Value = 8500
$ws.Range("N19").Value = -9076
$ws.Range("H22").Value = 1836.3334
$ws.Range("I22").Value = 1000
$ws.Range("J22").Value = 2254.5
$ws.Range("K22").Value = 1000
$ws.Range("L22").Value = 2254.5
$ws.Range("M22").Value = -471
$ws.Range("N22").Value = -3312.5
$ws.Range("H23").Value = 9000
$ws.Range("J23").Value = 9000
$ws.Range("L23").Value = 9000
$ws.Range("N23").Value = -9446
$ws.Range("H24").Value = 9943.846
$ws.Range("J24").Value = 9943.846
$ws.Range("L24").Value = 9943.846
$ws.Range("N24").Value = -10289.846
$ws.Range("H25").Value = 3900
$ws.Range("J25").Value = 3900
$ws.Range("L25").Value = 3900
$ws.Range("N25").Value = -4958
$ws.Range("H70").Value = 15644.333
$ws.Range("I70").Value = 102004
$ws.Range("J70").Value = 4849.375
$ws.Range("K70").Value = 102004
$ws.Range("L70").Value = 4849.375
$ws.Range("M70").Value = -101734
$ws.Range("N70").Value = -5389.375
$ws.Range("H73").Value = 15644.333
$ws.Range("I73").Value = 102004
$ws.Range("J73").Value = 4849.375
$ws.Range("K73").Value = 102004
$ws.Range("L73").Value = 4849.375
$ws.Range("M73").Value = -101068
$ws.Range("N73").Value = -6721.375
$ws.Range("H102").Value = 2517.9534
$ws.Range("I102").Value = 2652.9744
$ws.Range("J102").Value = 1201.5
$ws.Range("K102").Value = 2652.9744
$ws.Range("L102").Value = 1201.5
$ws.Range("M102").Value = -1030.9744
$ws.Range("N102").Value = -4445.5
$ws.Range("H126").Value = 3367.1333
$ws.Range("I126").Value = 2512.5
$ws.Range("J126").Value = 3936.889
$ws.Range("K126").Value = 7537.5
$ws.Range("L126").Value = 11810.667
$ws.Range("M126").Value = -5067.5
$ws.Range("N126").Value = -16750.667

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 6124.4287
$ws.Range("I7").Value = 7367.4443
$ws.Range("J7").Value = 5535.6313
$ws.Range("K7").Value = 7367.4443
$ws.Range("L7").Value = 5535.6313
$ws.Range("M7").Value = -7255.4443
$ws.Range("N7").Value = -5759.6313
$ws.Range("H40").Value = 5126.7393
$ws.Range("I40").Value = 5226.0625
$ws.Range("J40").Value = 4899.7144
$ws.Range("K40").Value = 5226.0625
$ws.Range("L40").Value = 4899.7144
$ws.Range("M40").Value = -5090.0625
$ws.Range("N40").Value = -5171.7144
$ws.Range("H122").Value = 6396
$ws.Range("I122").Value = 6418.5
$ws.Range("K122").Value = 19255.5
$ws.Range("M122").Value = -16805.5
$ws.Range("H126").Value = 6124.4287
$ws.Range("I126").Value = 7367.4443
$ws.Range("J126").Value = 5535.6313
$ws.Range("K126").Value = 22102.3329
$ws.Range("L126").Value = 16606.8939
$ws.Range("M126").Value = -19632.3329
$ws.Range("N126").Value = -21546.8939

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H15").Value = 5004.2
$ws.Range("J15").Value = 5004.2
$ws.Range("L15").Value = 5004.2
$ws.Range("N15").Value = -5580.2
$ws.Range("H17").Value = 985
$ws.Range("I17").Value = 975
$ws.Range("J17").Value = 1005
$ws.Range("K17").Value = 975
$ws.Range("L17").Value = 1005
$ws.Range("M17").Value = -803
$ws.Range("N17").Value = -1349
$ws.Range("H21").Value = 11000
$ws.Range("J21").Value = 11000
$ws.Range("L21").Value = 11000
$ws.Range("N21").Value = -11470
$ws.Range("H23").Value = 1381.5385
$ws.Range("I23").Value = 1177.1818
$ws.Range("J23").Value = 2505.5
$ws.Range("K23").Value = 1177.1818
$ws.Range("L23").Value = 2505.5
$ws.Range("M23").Value = -948.1818000000001
$ws.Range("N23").Value = -2963.5
$ws.Range("H24").Value = 6104.5
$ws.Range("I24").Value = 1209
$ws.Range("J24").Value = 11000
$ws.Range("K24").Value = 1209
$ws.Range("L24").Value = 11000
$ws.Range("M24").Value = -979
$ws.Range("N24").Value = -11460
$ws.Range("H35").Value = 11000
$ws.Range("J35").Value = 11000
$ws.Range("L35").Value = 11000
$ws.Range("N35").Value = -11580
$ws.Range("H126").Value = 3864.8572
$ws.Range("I126").Value = 2800.3
$ws.Range("J126").Value = 6526.25
$ws.Range("K126").Value = 8400.900000000001
$ws.Range("L126").Value = 19578.75
$ws.Range("M126").Value = -5930.900000000001
$ws.Range("N126").Value = -24518.75
